$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.043350035268205
$ws.Range("D2").Value = 1.041332631022898
$ws.Range("E2").Value = 1.056681566927112
$ws.Range("F2").Value = 1.06359659467957
$ws.Range("I2").Value = 1.033771100570638
$ws.Range("J2").Value = 1.048420956217138
$ws.Range("K2").Value = 1.044112444682736
$ws.Range("L2").Value = 1.059418610081964
$ws.Range("M2").Value = 1.066314805497059
$ws.Range("N2").Value = 1.019874624930268

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.044784002071997
$ws.Range("D3").Value = 1.042354381503289
$ws.Range("E3").Value = 1.058143109363738
$ws.Range("F3").Value = 1.065204312126877
$ws.Range("I3").Value = 1.034023081233812
$ws.Range("J3").Value = 1.049499332053239
$ws.Range("K3").Value = 1.044944246879556
$ws.Range("L3").Value = 1.060692191973851
$ws.Range("M3").Value = 1.067735577024751
$ws.Range("N3").Value = 1.020244682089662

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045710604599175
$ws.Range("D4").Value = 1.043014199200374
$ws.Range("E4").Value = 1.059087945585905
$ws.Range("F4").Value = 1.066243953244835
$ws.Range("I4").Value = 1.034184216774504
$ws.Range("J4").Value = 1.05019543228345
$ws.Range("K4").Value = 1.045480554169492
$ws.Range("L4").Value = 1.061514891988344
$ws.Range("M4").Value = 1.06865376559317
$ws.Range("N4").Value = 1.020483265786914

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046099850889578
$ws.Range("D5").Value = 1.043291273414001
$ws.Range("E5").Value = 1.059484950423417
$ws.Range("F5").Value = 1.066680867532015
$ws.Range("I5").Value = 1.034251501469906
$ws.Range("J5").Value = 1.050487675169966
$ws.Range("K5").Value = 1.045705560548673
$ws.Range("L5").Value = 1.06186042731162
$ws.Range("M5").Value = 1.069039504130961
$ws.Range("N5").Value = 1.020583360073455

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046165189805298
$ws.Range("D6").Value = 1.043337777135529
$ws.Range("E6").Value = 1.059551597482288
$ws.Range("F6").Value = 1.066754218722753
$ws.Range("I6").Value = 1.03426277212194
$ws.Range("J6").Value = 1.050536720856968
$ws.Range("K6").Value = 1.045743313399983
$ws.Range("L6").Value = 1.06191842515587
$ws.Range("M6").Value = 1.069104255800098
$ws.Range("N6").Value = 1.020600154290713

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045715806887614
$ws.Range("D7").Value = 1.043017902705559
$ws.Range("E7").Value = 1.059093251179055
$ws.Range("F7").Value = 1.066249791898866
$ws.Range("I7").Value = 1.034185117628837
$ws.Range("J7").Value = 1.050199338803776
$ws.Range("K7").Value = 1.045483562507954
$ws.Range("L7").Value = 1.061519510329161
$ws.Range("M7").Value = 1.068658920895371
$ws.Range("N7").Value = 1.020484604060425

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.043834917983032
$ws.Range("D8").Value = 1.041678212547009
$ws.Range("E8").Value = 1.05717568789799
$ws.Range("F8").Value = 1.064140071150009
$ws.Range("I8").Value = 1.033856655423532
$ws.Range("J8").Value = 1.04878574934566
$ws.Range("K8").Value = 1.044393956686222
$ws.Range("L8").Value = 1.059849314799953
$ws.Range("M8").Value = 1.066795202910922
$ws.Range("N8").Value = 1.019999868122309

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040510529101135
$ws.Range("D9").Value = 1.039307209723246
$ws.Range("E9").Value = 1.053789680512645
$ws.Range("F9").Value = 1.060417104731024
$ws.Range("I9").Value = 1.033263157161328
$ws.Range("J9").Value = 1.04628173360309
$ws.Range("K9").Value = 1.042459036567018
$ws.Range("L9").Value = 1.056895297610551
$ws.Range("M9").Value = 1.06350203539478
$ws.Range("N9").Value = 1.019138985379948

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03828711480803
$ws.Range("D10").Value = 1.037719386250764
$ws.Range("E10").Value = 1.051527231204207
$ws.Range("F10").Value = 1.057931071669741
$ws.Range("I10").Value = 1.032857524415149
$ws.Range("J10").Value = 1.04460328471994
$ws.Range("K10").Value = 1.041158852887813
$ws.Range("L10").Value = 1.05491825492642
$ws.Range("M10").Value = 1.061300100574099
$ws.Range("N10").Value = 1.018560451859899

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037322561509664
$ws.Range("D11").Value = 1.037030091798374
$ws.Range("E11").Value = 1.050546265103486
$ws.Range("F11").Value = 1.056853530254003
$ws.Range("I11").Value = 1.032679498679429
$ws.Range("J11").Value = 1.043874269998979
$ws.Range("K11").Value = 1.040593380827931
$ws.Range("L11").Value = 1.054060267845481
$ws.Range("M11").Value = 1.0603450096844
$ws.Range("N11").Value = 1.018308824897793

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.036964004886
$ws.Range("D12").Value = 1.036773788661865
$ws.Range("E12").Value = 1.050181686022158
$ws.Range("F12").Value = 1.056453113565733
$ws.Range("I12").Value = 1.032613012183325
$ws.Range("J12").Value = 1.043603139681706
$ws.Range("K12").Value = 1.040382961732405
$ws.Range("L12").Value = 1.053741278630382
$ws.Range("M12").Value = 1.059989992431925
$ws.Range("N12").Value = 1.018215189538723

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037040929236287
$ws.Range("D13").Value = 1.036828778778623
$ws.Range("E13").Value = 1.050259898855049
$ws.Range("F13").Value = 1.056539012129577
$ws.Range("I13").Value = 1.032627290070986
$ws.Range("J13").Value = 1.04366131361814
$ws.Range("K13").Value = 1.040428114500842
$ws.Range("L13").Value = 1.053809716344127
$ws.Range("M13").Value = 1.060066156456137
$ws.Range("N13").Value = 1.018235282360435

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037292928837446
$ws.Range("D14").Value = 1.037008911204977
$ws.Range("E14").Value = 1.050516133101222
$ws.Range("F14").Value = 1.056820435200464
$ws.Range("I14").Value = 1.032674010234446
$ws.Range("J14").Value = 1.043851865289468
$ws.Range("K14").Value = 1.040575995248319
$ws.Range("L14").Value = 1.054033906147947
$ws.Range("M14").Value = 1.060315669073764
$ws.Range("N14").Value = 1.018301088446884

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037448156950141
$ws.Range("D15").Value = 1.037119861022499
$ws.Range("E15").Value = 1.050673980135438
$ws.Range("F15").Value = 1.056993806478206
$ws.Range("I15").Value = 1.0327027483416
$ws.Range("J15").Value = 1.043969224993374
$ws.Range("K15").Value = 1.040667059289975
$ws.Range("L15").Value = 1.054171997636971
$ws.Range("M15").Value = 1.060469368144792
$ws.Range("N15").Value = 1.01834161117742

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038351090413365
$ws.Range("D16").Value = 1.03776509505573
$ws.Range("E16").Value = 1.05159230642538
$ws.Range("F16").Value = 1.058002561205755
$ws.Range("I16").Value = 1.032869289028292
$ws.Range("J16").Value = 1.04465161942315
$ws.Range("K16").Value = 1.041196328696625
$ws.Range("L16").Value = 1.054975155771979
$ws.Range("M16").Value = 1.061363451624079
$ws.Range("N16").Value = 1.01857712779185

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038916989032728
$ws.Range("D17").Value = 1.038169359838133
$ws.Range("E17").Value = 1.052167991944533
$ws.Range("F17").Value = 1.058635033103585
$ws.Range("I17").Value = 1.032973116168956
$ws.Range("J17").Value = 1.045079064775489
$ws.Range("K17").Value = 1.041527657287251
$ws.Range("L17").Value = 1.055478438165001
$ws.Range("M17").Value = 1.061923842390553
$ws.Range("N17").Value = 1.018724560439628

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039246894978223
$ws.Range("D18").Value = 1.038404991491892
$ws.Range("E18").Value = 1.052503653635058
$ws.Range("F18").Value = 1.059003840350868
$ws.Range("I18").Value = 1.033033446865614
$ws.Range("J18").Value = 1.045328171151652
$ws.Range("K18").Value = 1.041720676284113
$ws.Range("L18").Value = 1.055771810017669
$ws.Range("M18").Value = 1.062250551274075
$ws.Range("N18").Value = 1.018810447655902

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039359355256112
$ws.Range("D19").Value = 1.038485307243449
$ws.Range("E19").Value = 1.05261808446827
$ws.Range("F19").Value = 1.059129576897032
$ws.Range("I19").Value = 1.033053979124975
$ws.Range("J19").Value = 1.045413073655084
$ws.Range("K19").Value = 1.04178645030621
$ws.Range("L19").Value = 1.055871811195473
$ws.Range("M19").Value = 1.062361924068516
$ws.Range("N19").Value = 1.018839714757945

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.03885629143453
$ws.Range("D20").Value = 1.038126003582755
$ws.Range("E20").Value = 1.052106239418163
$ws.Range("F20").Value = 1.058567185561945
$ws.Range("I20").Value = 1.0329620002927
$ws.Range("J20").Value = 1.045033226243734
$ws.Range("K20").Value = 1.041492133655728
$ws.Range("L20").Value = 1.055424459835894
$ws.Range("M20").Value = 1.06186373411946
$ws.Range("N20").Value = 1.01870875347109

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037218728991528
$ws.Range("D21").Value = 1.036955874143814
$ws.Range("E21").Value = 1.050440684212014
$ws.Range("F21").Value = 1.056737567878409
$ws.Range("I21").Value = 1.032660262256496
$ws.Range("J21").Value = 1.043795762054862
$ws.Range("K21").Value = 1.040532458513083
$ws.Range("L21").Value = 1.053967896027467
$ws.Range("M21").Value = 1.060242200900503
$ws.Range("N21").Value = 1.01828171490018

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036187512637309
$ws.Range("D22").Value = 1.036218611634954
$ws.Range("E22").Value = 1.049392295981357
$ws.Range("F22").Value = 1.055586227706925
$ws.Range("I22").Value = 1.032468465001626
$ws.Range("J22").Value = 1.043015738721999
$ws.Range("K22").Value = 1.039926885802792
$ws.Range("L22").Value = 1.053050388529027
$ws.Range("M22").Value = 1.059221205310886
$ws.Range("N22").Value = 1.018012234842096

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036734335826211
$ws.Range("D23").Value = 1.036609597449143
$ws.Range("E23").Value = 1.049948181566871
$ws.Range("F23").Value = 1.056196671373599
$ws.Range("I23").Value = 1.032570338307432
$ws.Range("J23").Value = 1.043429433578952
$ws.Range("K23").Value = 1.04024812012471
$ws.Range("L23").Value = 1.053536940854725
$ws.Range("M23").Value = 1.059762596657765
$ws.Range("N23").Value = 1.018155185231155

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.038883718601725
$ws.Range("D24").Value = 1.038145594931281
$ws.Range("E24").Value = 1.05213414311614
$ws.Range("F24").Value = 1.058597843263826
$ws.Range("I24").Value = 1.032967023789239
$ws.Range("J24").Value = 1.045053939368534
$ws.Range("K24").Value = 1.041508185995563
$ws.Range("L24").Value = 1.055448850888936
$ws.Range("M24").Value = 1.061890894946988
$ws.Range("N24").Value = 1.018715896292823

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.04137119372864
$ws.Range("D25").Value = 1.039921414801505
$ws.Range("E25").Value = 1.054665916397533
$ws.Range("F25").Value = 1.061380262739911
$ws.Range("I25").Value = 1.033418341581835
$ws.Range("J25").Value = 1.046930664119579
$ws.Range("K25").Value = 1.042961047733418
$ws.Range("L25").Value = 1.05766031238979
$ws.Range("M25").Value = 1.064354512556675
$ws.Range("N25").Value = 1.019362349944304
